# "added 4wk low sales check" - refreshed forecast figures after adding a
# check for low 4-week sales, which shifted the per-week MyForecast values
# (and the dependent Inventory Coverage / Seasonality Index columns) on the
# "Forecast Comparison" sheet, and the aggregate figures on "Summary".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# New MyForecast (D), Inventory Coverage (H), Seasonality Index (L) values,
# keyed by worksheet row (row 2 = week W10 ... row 17 = week W25).
$data = @(
    @{ Row = 2;  D = 127; H = 10.28;              L = 0.91 },
    @{ Row = 3;  D = 127; H = 9.279999999999999;  L = 0.86 },
    @{ Row = 4;  D = 125; H = 8.41;               L = 1.12 },
    @{ Row = 5;  D = 124; H = 7.47;               L = 0.89 },
    @{ Row = 6;  D = 123; H = 6.52;               L = 0.84 },
    @{ Row = 7;  D = 122; H = 5.57;               L = 1.05 },
    @{ Row = 8;  D = 123; H = 4.53;               L = 1.16 },
    @{ Row = 9;  D = 122; H = 3.56;               L = 1.13 },
    @{ Row = 10; D = 121; H = 2.58;               L = 0.83 },
    @{ Row = 11; D = 120; H = 1.59;               L = 0.99 },
    @{ Row = 12; D = 118; H = 0.6;                L = 1.13 },
    @{ Row = 13; D = 118; H = $null;              L = 0.95 },
    @{ Row = 14; D = 117; H = $null;              L = 0.82 },
    @{ Row = 15; D = 117; H = $null;              L = 1.08 },
    @{ Row = 16; D = 115; H = $null;              L = 0.84 },
    @{ Row = 17; D = 115; H = $null;              L = 0.92 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    if ($item.H -ne $null) {
        $ws.Cells.Item($item.Row, 8).Value = $item.H
    }
    $ws.Cells.Item($item.Row, 12).Value = $item.L
}

# Update the "Summary" sheet aggregate figures, which are stored as text.
# Force text formatting so the written values stay strings (matching the
# original cell type) rather than being auto-coerced to numbers.
$summaryCells = @("B9", "B10", "B11", "B12", "B14")
foreach ($addr in $summaryCells) {
    $summary.Range($addr).NumberFormat = "@"
}

$summary.Range("B9").Value = "1934"
$summary.Range("B10").Value = "993"
$summary.Range("B11").Value = "503"
$summary.Range("B12").Value = "127"
$summary.Range("B14").Value = "115"
